$d = $word.ActiveDocument

# --- 1. Move the "_GoBack" bookmark ------------------------------------
# It currently sits right before the run "to understand this we have to
# learn all about synchronous and asynchronous programming." (i.e. right
# after "... In order to understand this ..."). Remove it from there -
# it gets re-created as an empty bookmark at the very end of the
# document further down below.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- 2. Drop the paragraph holding the embedded OLE/ActiveX control ----
# The last two paragraphs of the document are:
#   second-to-last : a paragraph whose only content is the embedded
#                    OLE/ActiveX control (little icon) that used to
#                    follow "... let's learn all about asynchronous
#                    programming."
#   last           : an empty trailing paragraph
# Removing the control paragraph merges it away, leaving a single empty
# paragraph at the end of the document.
$paraCount = $d.Paragraphs.Count
$controlPara = $d.Paragraphs.Item($paraCount - 1)
$controlPara.Range.Delete()

# --- 3. Re-insert the "_GoBack" bookmark in the final empty paragraph --
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
